$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (739) down across the
# new rows (740:751) so styles (date number format, fonts, etc.) match the
# rest of the table.
$ws.Range("A739:I739").Copy()
$ws.Range("A740:I751").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 739's "Localisation douleur" cell (G) holds text, so its style isn't
# right for rows whose G cell stays empty. Re-copy the blank-G style (taken
# from row 736, which has an empty G cell) onto those specific rows.
$ws.Range("G736").Copy()
$ws.Range("G741").PasteSpecial(-4122)
$ws.Range("G747:G750").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$nbsp = [char]0x00A0

# New training-log entries for the 2026-01-20 session (date serial 46042)
$data = @(
    @(740, "Maé Clavel",       6, 5, 3, "Ischio$nbsp",      6),
    @(741, "Naim Ighbane",     8, 6, 0, "",                 7),
    @(742, "Amir Etien",       7, 7, 3, "Ischio",           8),
    @(743, "Theo Owono",       7, 5, 3, "Coup cheville",   10),
    @(744, "Kamal Bafounta",   8, 4, 1, "Genou",           10),
    @(745, "Romain Thunet",    7, 5, 5, "Quadri",           3),
    @(746, "Omar Benyounes",   7, 3, 6, "Quadri",           6),
    @(747, "Yoan Zouma",       8, 8, 0, "",                10),
    @(748, "Jeremie Laurent",  8, 7, 0, "",                10),
    @(749, "Mattheo Haon",     9, 7, 0, "",                 9),
    @(750, "Ilan Ihaddadene",  9, 8, 0, "",                10),
    @(751, "Hedi Nasri",       8, 6, 2, "Hanche",           9)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = 46042
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = 70
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    if ($row[5] -ne "") {
        $ws.Cells.Item($r, 7).Value = $row[5]
    } else {
        $ws.Cells.Item($r, 7).Value = ""
    }
    $ws.Cells.Item($r, 8).Value = $row[6]
}

# Fill in the C*D formula down through the new rows (kept as its own
# shared-formula group so it doesn't disturb the existing I708:I739 group).
$ws.Range("I740:I751").Formula = "=C740*D740"

# Restore the workbook view state recorded in the saved file.
[void]$ws.Range("C755").Select()
